$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" '28.372.73'
Set-TextValue "E2" '  +3.38%  '

# Row 3
Set-TextValue "D3" '1.868.74'
Set-TextValue "E3" '  +1.75%  '

# Row 4
Set-TextValue "E4" '  -0.25%  '

# Row 5
Set-TextValue "D5" '339.05'

# Row 6
Set-TextValue "E6" '  -0.17%  '

# Row 7
Set-TextValue "D7" '0.4700'
Set-TextValue "E7" '  +2.03%  '

# Row 8
Set-TextValue "D8" '0.3968'
Set-TextValue "E8" '  +3.88%  '

# Row 9
Set-TextValue "D9" '47.57'
Set-TextValue "E9" '  +2.49%  '

# Row 10
Set-TextValue "D10" '0.08031'
Set-TextValue "E10" '  +1.60%  '

# Row 11
Set-TextValue "D11" '0.9998'
Set-TextValue "E11" '  +2.59%  '

# Row 12
Set-TextValue "D12" '21.94'
Set-TextValue "E12" '  +3.73%  '

# Row 13
Set-TextValue "D13" '6.027'
Set-TextValue "E13" '  +1.98%  '

# Row 14
Set-TextValue "D14" '1.868.31'
Set-TextValue "E14" '  -0.18%  '

# Row 15
Set-TextValue "D15" '7.246'
Set-TextValue "E15" '  +2.82%  '

# Row 16
Set-TextValue "D16" '91.26'
Set-TextValue "E16" '  +3.75%  '

# Row 17
Set-TextValue "D17" '1.003'
Set-TextValue "E17" '  +0.07%  '

# Row 18
Set-TextValue "E18" '  +1.38%  '

# Row 19
Set-TextValue "D19" '0.06619'
Set-TextValue "E19" '  +0.06%  '

# Row 20
Set-TextValue "D20" '17.55'
Set-TextValue "E20" '  +2.95%  '

# Row 22
Set-TextValue "D22" '28.380.72'
Set-TextValue "E22" '  +3.42%  '

# Row 23
Set-TextValue "D23" '5.476'
Set-TextValue "E23" '  +2.07%  '

# Row 24
Set-TextValue "E24" '  +1.81%  '

# Row 25
Set-TextValue "D25" '2.257'

# Row 26
Set-TextValue "D26" '2.089.08'
Set-TextValue "E26" '  +0.30%  '

# Row 27
Set-TextValue "D27" '161.22'
Set-TextValue "E27" '  +2.59%  '

# Row 28
Set-TextValue "D28" '19.75'
Set-TextValue "E28" '  +1.76%  '

# Row 29
Set-TextValue "D29" '2.121'
Set-TextValue "E29" '  +2.46%  '

# Row 30
Set-TextValue "D30" '5.501'
Set-TextValue "E30" '  +3.16%  '

# Row 31
Set-TextValue "D31" '120.40'
Set-TextValue "E31" '  +1.21%  '

# Row 32
Set-TextValue "D32" '0.9702'
Set-TextValue "E32" '  +1.54%  '

# Row 33
Set-TextValue "D33" '0.09508'
Set-TextValue "E33" '  +2.32%  '

# Row 34
Set-TextValue "D34" '3.598'
Set-TextValue "E34" '  +0.94%  '

# Row 35
Set-TextValue "B35" 'Filecoin'
Set-TextValue "C35" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D35" '5.351'
Set-TextValue "E35" '  +1.97%  '

# Row 36
Set-TextValue "B36" 'ARBITRUM'
Set-TextValue "C36" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D36" '1.373'
Set-TextValue "E36" '  +4.21%  '

# Row 37
Set-TextValue "E37" '  +2.66%  '

# Row 38
Set-TextValue "E38" '  +2.76%  '

# Row 39
Set-TextValue "D39" '8.384'
Set-TextValue "E39" '  +3.84%  '

# Row 40
Set-TextValue "E40" '  +2.34%  '

# Row 41
Set-TextValue "D41" '0.5955'
Set-TextValue "E41" '  +2.60%  '

# Row 42
Set-TextValue "E42" '  -0.15%  '

# Row 43
Set-TextValue "D43" '0.1873'
Set-TextValue "E43" '  +1.64%  '

# Row 44
Set-TextValue "D44" '10.36'
Set-TextValue "E44" '  +3.42%  '

# Row 45
Set-TextValue "D45" '1.290'
Set-TextValue "E45" '  +2.29%  '

# Row 46
Set-TextValue "D46" '0.5593'
Set-TextValue "E46" '  +1.82%  '

# Row 47
Set-TextValue "D47" '12.16'
Set-TextValue "E47" '  +1.64%  '

# Row 48
Set-TextValue "D48" '1.957'
Set-TextValue "E48" '  +4.63%  '

# Row 49
Set-TextValue "D49" '0.06875'
Set-TextValue "E49" '  +3.32%  '

# Row 50
Set-TextValue "D50" '2.083'
Set-TextValue "E50" '  +17.74%  '

# Row 51
Set-TextValue "D51" '111.68'
Set-TextValue "E51" '  +1.49%  '
